{"js": "// Update the worksheet's \"three-digit \u00f7 one-digit\" problems.\n// Each cell holds a single run/paragraph with a unique \"NNN\u00f7N=\" label, so we\n// can safely do a targeted search-and-replace per label, in document order,\n// using exact (case-sensitive, non-wildcard) matching.\nconst replacements = [\n  [\"450\u00f73=\", \"849\u00f75=\"],\n  [\"437\u00f79=\", \"348\u00f79=\"],\n  [\"595\u00f74=\", \"945\u00f75=\"],\n  [\"221\u00f73=\", \"229\u00f79=\"],\n  [\"679\u00f74=\", \"810\u00f76=\"],\n  [\"955\u00f75=\", \"566\u00f76=\"],\n  [\"113\u00f78=\", \"445\u00f78=\"],\n  [\"917\u00f78=\", \"125\u00f77=\"],\n  [\"106\u00f77=\", \"563\u00f74=\"],\n  [\"519\u00f73=\", \"151\u00f74=\"],\n  [\"958\u00f72=\", \"199\u00f79=\"],\n  [\"433\u00f73=\", \"226\u00f79=\"],\n  [\"753\u00f75=\", \"113\u00f78=\"],\n  [\"500\u00f79=\", \"977\u00f72=\"],\n  [\"766\u00f74=\", \"392\u00f79=\"],\n  [\"937\u00f74=\", \"220\u00f75=\"],\n  [\"175\u00f73=\", \"963\u00f78=\"],\n  [\"769\u00f76=\", \"741\u00f78=\"],\n  [\"119\u00f75=\", \"639\u00f79=\"],\n  [\"975\u00f79=\", \"494\u00f74=\"],\n  [\"913\u00f77=\", \"776\u00f77=\"],\n  [\"713\u00f79=\", \"562\u00f74=\"],\n  [\"915\u00f76=\", \"119\u00f79=\"],\n  [\"333\u00f77=\", \"768\u00f73=\"],\n  [\"292\u00f72=\", \"105\u00f75=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text to replace: \"${oldText}\"`);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the worksheet's \"three-digit \u00f7 one-digit\" problems.\n# Each cell holds a single, unique \"NNN\u00f7N=\" label, so a simple\n# Find/Replace (one match, exact text) per pair is sufficient and safe\n# even though a couple of the *new* values happen to equal an *older*\n# cell's original text -- by the time we get to those pairs the earlier\n# occurrence has already been rewritten, so Find still matches exactly\n# one cell.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"450\u00f73=\", \"849\u00f75=\"),\n    @(\"437\u00f79=\", \"348\u00f79=\"),\n    @(\"595\u00f74=\", \"945\u00f75=\"),\n    @(\"221\u00f73=\", \"229\u00f79=\"),\n    @(\"679\u00f74=\", \"810\u00f76=\"),\n    @(\"955\u00f75=\", \"566\u00f76=\"),\n    @(\"113\u00f78=\", \"445\u00f78=\"),\n    @(\"917\u00f78=\", \"125\u00f77=\"),\n    @(\"106\u00f77=\", \"563\u00f74=\"),\n    @(\"519\u00f73=\", \"151\u00f74=\"),\n    @(\"958\u00f72=\", \"199\u00f79=\"),\n    @(\"433\u00f73=\", \"226\u00f79=\"),\n    @(\"753\u00f75=\", \"113\u00f78=\"),\n    @(\"500\u00f79=\", \"977\u00f72=\"),\n    @(\"766\u00f74=\", \"392\u00f79=\"),\n    @(\"937\u00f74=\", \"220\u00f75=\"),\n    @(\"175\u00f73=\", \"963\u00f78=\"),\n    @(\"769\u00f76=\", \"741\u00f78=\"),\n    @(\"119\u00f75=\", \"639\u00f79=\"),\n    @(\"975\u00f79=\", \"494\u00f74=\"),\n    @(\"913\u00f77=\", \"776\u00f77=\"),\n    @(\"713\u00f79=\", \"562\u00f74=\"),\n    @(\"915\u00f76=\", \"119\u00f79=\"),\n    @(\"333\u00f77=\", \"768\u00f73=\"),\n    @(\"292\u00f72=\", \"105\u00f75=\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
